$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the production figures (columns C, N, P, R, V) for rows 4-13 ---
# Row 4
$ws.Range("C4").Value = 1827707
$ws.Range("N4").Value = 6.89
$ws.Range("P4").Value = 9.78
$ws.Range("R4").Value = 0.69
$ws.Range("V4").Value = 18

# Row 5
$ws.Range("C5").Value = 221717
$ws.Range("N5").Value = 6.34
$ws.Range("P5").Value = 9.78
$ws.Range("R5").Value = 0.67
$ws.Range("V5").Value = 18

# Row 6
$ws.Range("C6").Value = 2096181
$ws.Range("N6").Value = 6.28
$ws.Range("P6").Value = 10.67
$ws.Range("R6").Value = 0.724
$ws.Range("V6").Value = 18

# Row 7
$ws.Range("C7").Value = 1388642
$ws.Range("N7").Value = 9.1
$ws.Range("P7").Value = 8.94
$ws.Range("R7").Value = 0.61
$ws.Range("V7").Value = 18

# Row 8
$ws.Range("C8").Value = 1381714
$ws.Range("N8").Value = 8.41
$ws.Range("P8").Value = 9.78
$ws.Range("R8").Value = 0.681
$ws.Range("V8").Value = 18

# Row 9
$ws.Range("C9").Value = 1855634
$ws.Range("N9").Value = 7.24
$ws.Range("P9").Value = 9.78
$ws.Range("V9").Value = 18

# Row 10
$ws.Range("C10").Value = 2526024
$ws.Range("N10").Value = 9.49
$ws.Range("P10").Value = 9.78
$ws.Range("R10").Value = 0.654
$ws.Range("V10").Value = 18

# Row 11 (totals row - N/P/V are direct inputs, C/E/etc. are formulas)
$ws.Range("N11").Value = 7.82
$ws.Range("P11").Value = 9.78
$ws.Range("V11").Value = 18

# Row 12
$ws.Range("C12").Value = 2941602
$ws.Range("N12").Value = 3.71
$ws.Range("R12").Value = 0.81
$ws.Range("V12").Value = 18

# Row 13
$ws.Range("C13").Value = 523283
$ws.Range("P13").Value = 9.89
$ws.Range("V13").Value = 18

# --- Apply a number format to a few blank cells (matches the touched formatting
# around the lower summary table) ---
$ws.Range("M15:M16").NumberFormat = "#,##0"
$ws.Range("M18:M19").NumberFormat = "#,##0"
$ws.Range("J17").NumberFormat = "#,##0"

# --- Update the active selection to match the author's last cursor position ---
$ws.Range("C8").Select()
